# Ideas of features: Adding intensity-sensitive color vision characteristics
#
# This script reproduces the commit that inserts 5 new "lid" creature
# features (color-vision modifiers keyed to light intensity) into the
# "Feuil1" ideas sheet, right after the existing eye/vision rows, and
# records the corresponding gene -> feature-code lookups on the "Feuil2"
# combinatorics sheet.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Feuil1")
$ws2 = $wb.Worksheets.Item("Feuil2")

# Insert 5 fresh rows starting at row 18 (right after the "Pasture eyes"
# / VI8 row), pushing Moustaches, Pretty/Fancy, Lifetime/aging and the
# own-light features (Glowy/Shiny/Blushy) down by 5 rows.
$ws1.Rows("18:22").Insert()

# New row 18: Day blue lids (LI1)
$ws1.Range("A18").Value = "Day blue lids"
$ws1.Range("B18").Value = "LI1"
$ws1.Range("C18").Value = "Shutters high levels of blue during the day"
$ws1.Range("D18").Value = "Stops high level of blue in high intensity light"
$ws1.Range("E18").Value = "colorVision blue capped (upper) at 3 in high intensity light"
$ws1.Range("F18").Value = "CDDBC, DDDBC"

# New row 19: Night lids (LI2)
$ws1.Range("A19").Value = "Night lids"
$ws1.Range("B19").Value = "LI2"
$ws1.Range("C19").Value = "Implify low level colors during the night"
$ws1.Range("D19").Value = "Extend lower cap for all colors in low intensity light"
$ws1.Range("E19").Value = "colorVision += [(-1, 0), (-1, 0) (-1, 0)] in low intensity light"
$ws1.Range("F19").Value = "CBDDBC"

# New row 20: Stripped lids (LI3)
$ws1.Range("A20").Value = "Stripped lids"
$ws1.Range("B20").Value = "LI3"
$ws1.Range("C20").Value = "Mix up the visible colors during the day"
$ws1.Range("D20").Value = "Shift values for each color vision to the right (values for blue become for red, etc) in high intensity light"
$ws1.Range("E20").Value = "colorVision = [(tuple3), (tuple1), (tuple2)] in high intensity light"
$ws1.Range("F20").Value = "CDCBC"

# New row 21: Marron lids (LI4)
$ws1.Range("A21").Value = "Marron lids"
$ws1.Range("B21").Value = "LI4"
$ws1.Range("C21").Value = "Stops blue light in the night"
$ws1.Range("D21").Value = "Stops all blue in low intensity light"
$ws1.Range("E21").Value = "colorVision += [(0, 0), (0, 0), (-100, -100)] in low intensity light"
$ws1.Range("F21").Value = "CBBBBC"

# New row 22: Ultragreen lids (LI5)
$ws1.Range("A22").Value = "Ultragreen lids"
$ws1.Range("B22").Value = "LI5"
$ws1.Range("C22").Value = "Boosts green vision in the day"
$ws1.Range("D22").Value = "Enable vision for green in high intensity light"
$ws1.Range("E22").Value = "colorVision += [(0, 0), (-100, 100), (0, 0)] in high intensity light"
$ws1.Range("F22").Value = "BDDCB"

# Update the sheet1 selection to match the authored state.
$ws1.Activate()
$ws1.Range("E38").Select()

# Record the new gene -> feature-code lookups on the combinatorics sheet
# (Feuil2) for the newly introduced LI1-LI5 lid features.
$ws2.Range("K53").Value  = "LI5"
$ws2.Range("K99").Value  = "LI3"
$ws2.Range("K108").Value = "LI1"
$ws2.Range("K109").Value = "LI1"
$ws2.Range("M246").Value = "LI4"
$ws2.Range("M318").Value = "LI2"

# Update the sheet2 scroll position / selection to match the authored state.
$ws2.Activate()
$ws2.Application.ActiveWindow.ScrollRow = 231
$ws2.Range("K257").Select()
